$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "updated activity till excel form" - runs/balls figures for rows 3 and 4
# were swapped between each other.
# Force text storage (cells already hold these numeric-looking figures as
# text) so "1"/"7"/"10" aren't silently reinterpreted as numbers.
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "1"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "7"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "10"
